$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new, empty paragraph at the very start of the document, with
#    the same paragraph formatting (ind left=-567 hanging=567, rPr noProof)
#    that the paragraph which used to be first already carries.
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
[void]$firstPara.Range.InsertParagraphBefore()

# ---------------------------------------------------------------------------
# 2. "Acquista prodotto" used to be preceded by a stale
#    <w:lastRenderedPageBreak/>; re-typing the run's text regenerates the
#    run without that stale marker.
# ---------------------------------------------------------------------------
[void]$d.Content.Find.Execute("Acquista prodotto", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Acquista prodotto", 2)

# ---------------------------------------------------------------------------
# 3. "formCart, permette di visualizzare ..." used to be split across two
#    runs (", " and "permette di visualizzare..."); merge them into one.
# ---------------------------------------------------------------------------
[void]$d.Content.Find.Execute(", permette di visualizzare i prodotti presenti nel carrello", $true, $false, $false, $false, $false, `
    $true, 1, $false, ", permette di visualizzare i prodotti presenti nel carrello", 2)

# ---------------------------------------------------------------------------
# 4. "Checkout, form usato per ..." used to be split across two runs
#    (", " and "form usato per..."); merge them into one.
# ---------------------------------------------------------------------------
[void]$d.Content.Find.Execute(", form usato per l’input dell’ informazione per l’acquisto quando viene premuto il bottone checkout del formCart", $true, $false, $false, $false, $false, `
    $true, 1, $false, ", form usato per l’input dell’ informazione per l’acquisto quando viene premuto il bottone checkout del formCart", 2)

# ---------------------------------------------------------------------------
# 5. "riepilogoOrdine, contiene le informazioni relative all'ordine ..." used
#    to be split into three runs: ", contiene le informazioni ",
#    "relative all'ordine" and " un bottone per annullare...". Only the
#    first two must merge; the third run must stay separate. A plain text
#    replacement across the whole span would coalesce all three runs into
#    one (this engine merges every adjacent, equally-formatted run it
#    touches during an edit), so instead: cut the third run's text away,
#    merge the first two runs with a retype, then re-append the third run's
#    original text so it becomes a run of its own again.
# ---------------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.StartsWith("riepilogoOrdine")) {
        $targetPara = $candidate
        break
    }
}

$anchor = $targetPara.Range.Start
$paraEnd = $targetPara.Range.End - 1
$run4Start = $anchor + ("riepilogoOrdine" + ", contiene le informazioni relative all’ordine").Length

$run4 = $d.Range($run4Start, $paraEnd)
$run4Text = $run4.Text
$run4.Delete()

$mergeRange = $d.Range($anchor, $anchor + ("riepilogoOrdine" + ", contiene le informazioni relative all’ordine").Length)
$mergeRange.Text = $mergeRange.Text

$targetPara.Range.InsertAfter($run4Text)
